# Applies the tracker_resultados.xlsx update:
#  - fills in the pending result for row 70 (G70/H70)
#  - appends 7 newly-tracked matches as rows 82-88 (results pending)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 70: settle pending result ---
$ws.Range("G70").Value = "Fallo"
$ws.Range("H70").Value = -1

# --- New rows 82-88: newly tracked matches (result columns still pending) ---
# Row 82
$ws.Range("A82").Value = 14601372
$ws.Range("B82").Formula = "'2025-09-12"
$ws.Range("B82").Style = "Normal"
$ws.Range("C82").Formula = "'Vit Kopriva"
$ws.Range("C82").Style = "Normal"
$ws.Range("D82").Formula = "'Pablo Llamas Ruiz"
$ws.Range("D82").Style = "Normal"
$ws.Range("E82").Formula = "'Gana Vit Kopriva"
$ws.Range("E82").Style = "Normal"
$ws.Range("F82").Value = 1.8
$ws.Range("G82").Formula = "'"
$ws.Range("G82").Style = "Normal"
$ws.Range("H82").Formula = "'"
$ws.Range("H82").Style = "Normal"

# Row 83
$ws.Range("A83").Value = 14601410
$ws.Range("B83").Formula = "'2025-09-11"
$ws.Range("B83").Style = "Normal"
$ws.Range("C83").Formula = "'Cannon Kingsley"
$ws.Range("C83").Style = "Normal"
$ws.Range("D83").Formula = "'Titouan Droguet"
$ws.Range("D83").Style = "Normal"
$ws.Range("E83").Formula = "'Gana Cannon Kingsley"
$ws.Range("E83").Style = "Normal"
$ws.Range("F83").Value = 4.5
$ws.Range("G83").Formula = "'"
$ws.Range("G83").Style = "Normal"
$ws.Range("H83").Formula = "'"
$ws.Range("H83").Style = "Normal"

# Row 84
$ws.Range("A84").Value = 14601396
$ws.Range("B84").Formula = "'2025-09-12"
$ws.Range("B84").Style = "Normal"
$ws.Range("C84").Formula = "'Gianluca Cadenasso"
$ws.Range("C84").Style = "Normal"
$ws.Range("D84").Formula = "'Jelle Sels"
$ws.Range("D84").Style = "Normal"
$ws.Range("E84").Formula = "'Gana Gianluca Cadenasso"
$ws.Range("E84").Style = "Normal"
$ws.Range("F84").Value = 1.53
$ws.Range("G84").Formula = "'"
$ws.Range("G84").Style = "Normal"
$ws.Range("H84").Formula = "'"
$ws.Range("H84").Style = "Normal"

# Row 85
$ws.Range("A85").Value = 14601395
$ws.Range("B85").Formula = "'2025-09-12"
$ws.Range("B85").Style = "Normal"
$ws.Range("C85").Formula = "'Marko Topo"
$ws.Range("C85").Style = "Normal"
$ws.Range("D85").Formula = "'Kilian Feldbausch"
$ws.Range("D85").Style = "Normal"
$ws.Range("E85").Formula = "'Gana Marko Topo"
$ws.Range("E85").Style = "Normal"
$ws.Range("F85").Value = 2
$ws.Range("G85").Formula = "'"
$ws.Range("G85").Style = "Normal"
$ws.Range("H85").Formula = "'"
$ws.Range("H85").Style = "Normal"

# Row 86
$ws.Range("A86").Value = 14601438
$ws.Range("B86").Formula = "'2025-09-11"
$ws.Range("B86").Style = "Normal"
$ws.Range("C86").Formula = "'Alex Marti Pujolras"
$ws.Range("C86").Style = "Normal"
$ws.Range("D86").Formula = "'Stefan Palosi"
$ws.Range("D86").Style = "Normal"
$ws.Range("E86").Formula = "'Gana Stefan Palosi"
$ws.Range("E86").Style = "Normal"
$ws.Range("F86").Value = 2.38
$ws.Range("G86").Formula = "'"
$ws.Range("G86").Style = "Normal"
$ws.Range("H86").Formula = "'"
$ws.Range("H86").Style = "Normal"

# Row 87
$ws.Range("A87").Value = 14662293
$ws.Range("B87").Formula = "'2025-09-11"
$ws.Range("B87").Style = "Normal"
$ws.Range("C87").Formula = "'Daniel De Jonge"
$ws.Range("C87").Style = "Normal"
$ws.Range("D87").Formula = "'Joshua Peck"
$ws.Range("D87").Style = "Normal"
$ws.Range("E87").Formula = "'Gana Joshua Peck"
$ws.Range("E87").Style = "Normal"
$ws.Range("F87").Value = 4.33
$ws.Range("G87").Formula = "'"
$ws.Range("G87").Style = "Normal"
$ws.Range("H87").Formula = "'"
$ws.Range("H87").Style = "Normal"

# Row 88
$ws.Range("A88").Value = 14655208
$ws.Range("B88").Formula = "'2025-09-11"
$ws.Range("B88").Style = "Normal"
$ws.Range("C88").Formula = "'Igor Marcondes"
$ws.Range("C88").Style = "Normal"
$ws.Range("D88").Formula = "'Daniel Dutra Da Silva"
$ws.Range("D88").Style = "Normal"
$ws.Range("E88").Formula = "'Gana Igor Marcondes"
$ws.Range("E88").Style = "Normal"
$ws.Range("F88").Value = 2.25
$ws.Range("G88").Formula = "'"
$ws.Range("G88").Style = "Normal"
$ws.Range("H88").Formula = "'"
$ws.Range("H88").Style = "Normal"
